$d = $word.ActiveDocument
$cell = $d.Tables(3).Cell(5, 7)
$cell.Range.Text = "0"
